$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.033.61'
$ws.Range('E2').Value = '  -0.09%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.386.96'
$ws.Range('E3').Value = '  -1.25%  '
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '571.64'
$ws.Range('E5').Value = '  -0.27%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.92'
$ws.Range('E6').Value = '  -0.50%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  -0.62%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.59'
$ws.Range('E9').Value = '  +0.56%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.123'
$ws.Range('E10').Value = '  -1.76%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.393'
$ws.Range('E11').Value = '  +1.59%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '3.971.08'
$ws.Range('E12').Value = '  -1.33%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.125'
$ws.Range('E13').Value = '  +2.35%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.71'
$ws.Range('E14').Value = '  -1.85%  '
$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000170'
$ws.Range('E15').Value = '  -1.14%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.383.03'
$ws.Range('E16').Value = '  -2.10%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.125.72'
$ws.Range('E17').Value = '  -0.42%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.11'
$ws.Range('E18').Value = '  -4.24%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.73'
$ws.Range('E19').Value = '  -4.74%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.92'
$ws.Range('E20').Value = '  -5.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '380.87'
$ws.Range('E21').Value = '  -4.12%  '
$ws.Range('B22').Value = 'Litecoin'
$ws.Range('C22').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '74.70'
$ws.Range('E22').Value = '  +2.21%  '
$ws.Range('B23').Value = 'Polygon'
$ws.Range('C23').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.554'
$ws.Range('E23').Value = '  -1.91%  '
$ws.Range('E24').Value = '  +0.21%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000117'
$ws.Range('E25').Value = '  -4.88%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.531.03'
$ws.Range('E26').Value = '  -1.67%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.181'
$ws.Range('E27').Value = '  +1.46%  '
$ws.Range('E28').Value = '  -0.31%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.32'
$ws.Range('E29').Value = '  -2.97%  '
$ws.Range('E30').Value = '  -0.88%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.97'
$ws.Range('E31').Value = '  -2.14%  '
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.40'
$ws.Range('E32').Value = '  -3.15%  '
$ws.Range('B33').Value = 'USDe'
$ws.Range('C33').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('E33').Value = '  -0.08%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '23.41'
$ws.Range('E34').Value = '  -2.01%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.96'
$ws.Range('E35').Value = '  -0.61%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '166.52'
$ws.Range('E36').Value = '  -0.23%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.421.70'
$ws.Range('E37').Value = '  -1.18%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.00'
$ws.Range('E38').Value = '  -2.63%  '
$ws.Range('E39').Value = '  -4.82%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0769'
$ws.Range('E40').Value = '  -2.03%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '27.31'
$ws.Range('E41').Value = '  -1.56%  '
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  -0.51%  '
$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.780'
$ws.Range('E43').Value = '  -2.66%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.85'
$ws.Range('E44').Value = '  -0.90%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.39'
$ws.Range('E45').Value = '  -2.50%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.66'
$ws.Range('E46').Value = '  -3.30%  '
$ws.Range('E47').Value = '  -0.90%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.460.16'
$ws.Range('E48').Value = '  -5.70%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '22.94'
$ws.Range('E49').Value = '  +0.32%  '
$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.77'
$ws.Range('E50').Value = '  -2.70%  '
$ws.Range('E51').Value = '  +2.47%  '
